# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.081.78'
$ws.Range("E2").Value = '  -2.68%  '
$ws.Range("D3").Value = '1.732.31'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '''310.64'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '''0.4874'
$ws.Range("E7").Value = '  +6.82%  '
$ws.Range("D8").Value = '''0.3517'
$ws.Range("E8").Value = '  +0.80%  '
$ws.Range("D9").Value = '''43.34'
$ws.Range("E9").Value = '  +3.30%  '
$ws.Range("D10").Value = '''0.07293'
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("D11").Value = '''1.054'
$ws.Range("E11").Value = '  -2.67%  '
$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("E13").Value = '  -2.78%  '
$ws.Range("E14").Value = '  -1.27%  '
$ws.Range("D15").Value = '1.731.27'
$ws.Range("D16").Value = '''6.907'
$ws.Range("E16").Value = '  -3.65%  '
$ws.Range("D17").Value = '''87.43'
$ws.Range("E17").Value = '  -4.49%  '
$ws.Range("D18").Value = '''0.00001041'
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").Value = '''0.06414'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").Value = '''0.9998'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").Value = '''16.61'
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("D22").Value = '''5.708'
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").Value = '27.130.08'
$ws.Range("E23").Value = '  -2.63%  '
$ws.Range("D24").Value = '''10.92'
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("D25").Value = '''2.084'
$ws.Range("E25").Value = '  -3.51%  '
$ws.Range("D26").Value = '''154.38'
$ws.Range("E26").Value = '  -4.96%  '
$ws.Range("D27").Value = '''20.03'
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '1.928.72'
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("D29").Value = '''2.090'
$ws.Range("E29").Value = '  -3.24%  '
$ws.Range("D30").Value = '''121.72'
$ws.Range("E30").Value = '  -1.25%  '
$ws.Range("D31").Value = '''1.049'
$ws.Range("E31").Value = '  -2.86%  '
$ws.Range("D32").Value = '''0.09327'
$ws.Range("E32").Value = '  +0.59%  '
$ws.Range("D33").Value = '''3.655'
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("D34").Value = '''5.425'
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.05939'
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '''0.02196'
$ws.Range("E36").Value = '  -2.37%  '
$ws.Range("D37").Value = '''11.04'
$ws.Range("E37").Value = '  -5.90%  '
$ws.Range("D38").Value = '''1.432'
$ws.Range("E38").Value = '  +4.64%  '
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("D40").Value = '''0.2003'
$ws.Range("E40").Value = '  -2.74%  '
$ws.Range("D41").Value = '''0.6016'
$ws.Range("E41").Value = '  -2.74%  '
$ws.Range("D42").Value = '''0.9995'
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '''1.097'
$ws.Range("E43").Value = '  -6.85%  '
$ws.Range("D44").Value = '''7.500'
$ws.Range("D45").Value = '''12.87'
$ws.Range("E45").Value = '  -1.15%  '
$ws.Range("D46").Value = '''3.587'
$ws.Range("E46").Value = '  -3.75%  '
$ws.Range("D47").Value = '''0.5685'
$ws.Range("E47").Value = '  -1.99%  '
$ws.Range("D48").Value = '''118.86'
$ws.Range("E48").Value = '  -2.76%  '
$ws.Range("D49").Value = '''1.855'
$ws.Range("E49").Value = '  -3.64%  '
$ws.Range("D50").Value = '''1.107'
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("D51").Value = '''0.06653'
$ws.Range("E51").Value = '  -1.93%  '
